$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 4570
$ws.Range("G2").Value = 80
$ws.Range("G3").Value = "不可售"
$ws.Range("G5").Value = 85
$ws.Range("F6").Value = 1066
$ws.Range("F9").Value = 372
$ws.Range("F10").Value = 367
$ws.Range("F11").Value = 2553
$ws.Range("F12").Value = 1287
$ws.Range("F16").Value = 21
$ws.Range("F20").Value = 10558
$ws.Range("F21").Value = 6123
$ws.Range("F26").Value = 6
$ws.Range("F28").Value = 844
$ws.Range("F30").Value = 183
$ws.Range("F31").Value = 863
$ws.Range("F32").Value = 3567
$ws.Range("F36").Value = 132
$ws.Range("F37").Value = 276
$ws.Range("F38").Value = 248
$ws.Range("F39").Value = 254
$ws.Range("F40").Value = 4868
$ws.Range("F42").Value = 1151
$ws.Range("F43").Value = 169
$ws.Range("F44").Value = 198
$ws.Range("F45").Value = 111
$ws.Range("F46").Value = 494

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F15").Value = 3604

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 8845
$ws.Range("F4").Value = 1664

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F4").Value = 4570
$ws.Range("G4").Value = 80
$ws.Range("G7").Value = 85
$ws.Range("F10").Value = 2553
$ws.Range("F14").Value = 1287
$ws.Range("F16").Value = 21
$ws.Range("F21").Value = 10558
$ws.Range("F22").Value = 3604
$ws.Range("F28").Value = 6
$ws.Range("F30").Value = 844
$ws.Range("F32").Value = 183
$ws.Range("F33").Value = 863
$ws.Range("F34").Value = 3567
$ws.Range("F36").Value = 132
$ws.Range("F37").Value = 276
$ws.Range("F38").Value = 248
$ws.Range("F40").Value = 254
$ws.Range("F41").Value = 4868
$ws.Range("F43").Value = 1151
$ws.Range("F44").Value = 169
$ws.Range("F45").Value = 494
